$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrdenSalida")

for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 28)
    # Force the value to remain text (not get auto-converted to a number)
    $cell.NumberFormat = "@"
    $cell.Value = "3536"
}
